$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (existing): A2=3, B2=365
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = 365

# Insert new data: row3 = (0,259), row4 = (1,204), row5 = (2,173)
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = 259

$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = 204

$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 173

# Apply the same formatting as A2 (style index 1) to the new A column cells
$ws.Range("A2").Copy()
$ws.Range("A3:A5").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
